$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the run-split sentence in the "Viewing even a simple line
# graph..." paragraph into a single run: "utilized" + ". " + "I chose..." +
# "A four-year..." -> one run "utilized. I chose... product."
# Matching the substring starting right at the run boundary after "utilized"
# (".  I chose...") lets the engine coalesce that run together with the
# "utilized" run on its left, while leaving the preceding " " run (the one
# right before "utilized") untouched, exactly mirroring the target XML.
# ---------------------------------------------------------------------------
$old1 = ". I chose to look at the average sales over time based on a quarterly and monthly basis. A four-year time-lapse seemed to offer an optimal compromise between identifying both long-term and short-term patterns for the individual products and any possible relationships they may share. Implementing the selectable key from Tableau also allows for a closer look at a particular product."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: add a new paragraph after the "The yellow product..." paragraph
# with the "Many of these patterns..." text, split across nine runs (to
# mirror the source document's run layout exactly).
# ---------------------------------------------------------------------------
$chunks = @(
    "Many of these patterns and the possible logic that I have discussed here can have a",
    "n enormous",
    " impact ",
    "on the Analyst, helping determine what avenues to pursue",
    " to continue their evaluation. These theories are not fact and would require further an",
    "aly",
    "sis to ensure validi",
    "t",
    "y."
)

# Locate the "yellow product" paragraph by its text and insert a new, empty
# paragraph right after it.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*reasonably constant sales rate year-round*") {
        $anchor = $p
    }
}
$anchor.Range.InsertParagraphAfter()
$targetIndex = $anchor.Index + 1

# Write the first chunk directly into the freshly created paragraph.
$p = $d.Paragraphs($targetIndex)
$p.Range.InsertAfter($chunks[0])

# For every following chunk: split off a new paragraph, insert the chunk
# text into it, then delete the paragraph mark that separates the two
# paragraphs so they collapse back into a single paragraph while keeping
# each piece of text in its own run (no formatting residue).
for ($i = 1; $i -lt $chunks.Length; $i++) {
    $p = $d.Paragraphs($targetIndex)
    $endPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $endPoint.InsertParagraphAfter()

    $newP = $d.Paragraphs($targetIndex + 1)
    $newP.Range.InsertAfter($chunks[$i])

    $markStart = $d.Paragraphs($targetIndex).Range.End - 1
    $markRange = $d.Range($markStart, $markStart + 1)
    $markRange.Delete()
}
